$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "comp_quantity_inst1": add a "type" classification column (E) that
# labels each origin/destination pair with a letter A..I.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")

$wsComp.Range("E1").Value = "type"
$wsComp.Range("E2").Value = "A"
$wsComp.Range("E3").Value = "B"
$wsComp.Range("E4").Value = "C"
$wsComp.Range("E5").Value = "D"
$wsComp.Range("E6").Value = "E"
$wsComp.Range("E7").Value = "F"
$wsComp.Range("E8").Value = "G"
$wsComp.Range("E9").Value = "H"
$wsComp.Range("E10").Value = "I"

# Center-align the new column (including two trailing blank rows so the
# used range grows down to row 12, matching the new table footprint) ...
$wsComp.Range("E1:E12").HorizontalAlignment = -4108
# ... then center-align the original data table as well.
$wsComp.Range("A1:D10").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Sheet "parameters": refresh the fleet-size figure and add two roll-up rows
# that total the "required" (no_req_total) vs "optional" (no_opt_total)
# shipment quantities pulled from comp_quantity_inst1.
# ---------------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")

$wsParams.Range("B12").Value = 22

$wsParams.Range("A13").Value = "no_req_total"
$wsParams.Range("B13").Formula = "=SUM(comp_quantity_inst1!C2:C6)"

$wsParams.Range("A14").Value = "no_opt_total"
$wsParams.Range("B14").Formula = "=SUM(comp_quantity_inst1!C7:C10)"

# ---------------------------------------------------------------------------
# Window/selection housekeeping: "parameters" becomes the active sheet (with
# the cursor parked just below the new rows); comp_quantity_inst1 keeps its
# own last-used selection but is no longer the active tab.
# ---------------------------------------------------------------------------
$wsComp.Activate() | Out-Null
$wsComp.Range("C11").Select() | Out-Null

$wsParams.Activate() | Out-Null
$wsParams.Range("B15").Select() | Out-Null

$wb.Save()
